# Add a new "Pull Request 4" section of video tutorial entries,
# following the same layout/formatting pattern used for the
# existing "Pull Request 2"/"Pull Request 3" sections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: section header "Pull Request 4:" - copy the green-fill
# formatting used by the other section header in row 10.
$ws.Range("A10").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Pull Request 4:"

# Topic names first (column A), matching the order the author typed
# them in before filling in the Google Drive links.
$ws.Range("A19").Value = "gitignore command"
$ws.Range("A20").Value = "git diff command"

# Row 19: gitignore command link
$ws.Range("C19").Value = "https://drive.google.com/file/d/1hSmhbMRXrPMPOj8Q5qEysc0L53-3dP1J/view?usp=sharing"

# Row 20: git diff command link
$ws.Range("C20").Value = "https://drive.google.com/file/d/1MwV1Isweru4x1E1YFr6yL91HEI-Z_AL9/view?usp=sharing"

# Move the active selection to D19, matching the author's final
# cursor position when they saved the workbook.
$ws.Range("D19").Select()
